$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.741.02'
$ws.Range('E2').Value = '  +3.06%  '
$ws.Range('D3').Value = '3.127.78'
$ws.Range('E3').Value = '  +1.75%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.70'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.26'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.82%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '3.125.14'
$ws.Range('E8').Value = '  +1.97%  '
$ws.Range('E9').Value = '  +0.79%  '
$ws.Range('E10').Value = '  +15.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.70'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('E12').Value = '  +0.91%  '
$ws.Range('E13').Value = '  +5.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.20'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.75%  '
$ws.Range('E15').Value = '  -0.69%  '
$ws.Range('D16').Value = '3.644.75'
$ws.Range('E16').Value = '  +1.75%  '
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('D18').Value = '63.623.09'
$ws.Range('E18').Value = '  +2.99%  '
$ws.Range('D19').Value = '3.123.32'
$ws.Range('E19').Value = '  +1.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '465.77'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.35'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.65%  '
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('E23').Value = '  +1.41%  '
$ws.Range('E24').Value = '  -3.20%  '
$ws.Range('E25').Value = '  +0.83%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +9.97%  '
$ws.Range('E28').Value = '  +1.90%  '
$ws.Range('E29').Value = '  -1.39%  '
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.88'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.15'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.34%  '
$ws.Range('E33').Value = '  -4.48%  '
$ws.Range('E34').Value = '  +11.06%  '
$ws.Range('E35').Value = '  +8.63%  '
$ws.Range('E36').Value = '  +1.81%  '
$ws.Range('E37').Value = '  +16.71%  '
$ws.Range('E38').Value = '  +1.54%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '457.15'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.71%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '50.99'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.75'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.58%  '
$ws.Range('E42').Value = '  +1.49%  '
$ws.Range('D43').Value = '2.911.52'
$ws.Range('E43').Value = '  -1.11%  '
$ws.Range('E44').Value = '  +1.02%  '
$ws.Range('E45').Value = '  +2.18%  '
$ws.Range('E46').Value = '  +3.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '127.15'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.82%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '35.72'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.71%  '
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('E50').Value = '  +0.51%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.74'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.63%  '
